$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TestCell {
    param($addr, $value, $colorBgr, $forceText)
    $rng = $ws.Range($addr)
    if ($forceText) {
        $rng.Value = "'" + $value
    } else {
        $rng.Value = $value
    }
    $rng.WrapText = $true
    $rng.HorizontalAlignment = -4131
    $rng.VerticalAlignment = -4160
    $rng.Borders.LineStyle = 1
    $rng.Font.Color = $colorBgr
}

Set-TestCell "A9" "TC_03" 32768 $false
Set-TestCell "B9" "Verify the geopoliticalTypeNm record with parameter geopoliticalTypeNm with  geopoliticalTypeId, geopoliticalTypeName  attributes." 32768 $false
Set-TestCell "C9" "GraphQL" 32768 $false
$cellVal = @"

{
	"query":"
	{
		  geopoliticalTypes (geopoliticalTypeName : \"GeoTest\") 
		{
			    geopoliticalTypeId    
		}
	}"
}
"@
Set-TestCell "D9" $cellVal 32768 $false
$rng = $ws.Range("E9")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 32768
$rng = $ws.Range("F9")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 32768
$rng = $ws.Range("G9")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 32768
Set-TestCell "H9" "200" 32768 $true
$rng = $ws.Range("I9")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 32768
Set-TestCell "J9" "Pass" 32768 $false
Set-TestCell "K9" "Total number of records matching between DB & Response: 1, below are the test steps for this test case" 32768 $false
$rng = $ws.Range("A10")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 32768
Set-TestCell "B10" "Verify the geopoliticalTypeNm record with parameter geopoliticalTypeNm with  geopoliticalTypeId, geopoliticalTypeName  attributes." 32768 $false
Set-TestCell "C10" "GraphQL" 32768 $false
$cellVal = @"

{
	"query":"
	{
		  geopoliticalTypes (geopoliticalTypeName : \"GeoTest\") 
		{
			    geopoliticalTypeId    
		}
	}"
}
"@
Set-TestCell "D10" $cellVal 32768 $false
$rng = $ws.Range("E10")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 32768
$rng = $ws.Range("F10")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 32768
$rng = $ws.Range("G10")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 32768
$rng = $ws.Range("H10")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 32768
$cellVal = @"
Response_GeopoliticalType_ID: 5906024221803810540
DB_GeopoliticalType_ID: 5906024221803810540

"@
Set-TestCell "I10" $cellVal 32768 $false
Set-TestCell "J10" "Pass" 32768 $false
$rng = $ws.Range("K10")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 32768
Set-TestCell "A11" "TC_04" 255 $false
Set-TestCell "B11" "Verify no results fetched when passing the invalid geopoliticalTypeName  parameter." 255 $false
Set-TestCell "C11" "GraphQL" 255 $false
$rng = $ws.Range("D11")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 255
$rng = $ws.Range("E11")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 255
$rng = $ws.Range("F11")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 255
$rng = $ws.Range("G11")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 255
Set-TestCell "H11" "200" 255 $true
$cellVal = @"

{
	"meta":
	{
		"version":"1.0.0",
		"errors":
		[
			{
				"timestamp":"2020-02-04T13:38:34.837Z
				[
					GMT
				]",
				"error":"ValidationError",
				"message":"Validation error of type FieldUndefined: Field 'geopoliticalTypeId1' in type 'GeoPoliticalType' is undefined @ 'geopoliticalTypes/geopoliticalTypeId1'",
				"path":null
			}
		]
	},
	"data":null
}
"@
Set-TestCell "I11" $cellVal 255 $false
Set-TestCell "J11" "Fail" 255 $false
$rng = $ws.Range("K11")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 255
Set-TestCell "A12" "TC_05" 255 $false
Set-TestCell "B12" "Verify the error message when passing the multiple(2 attributes) invalid attribute names." 255 $false
Set-TestCell "C12" "GraphQL" 255 $false
$rng = $ws.Range("D12")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 255
$rng = $ws.Range("E12")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 255
$rng = $ws.Range("F12")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 255
$rng = $ws.Range("G12")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 255
Set-TestCell "H12" "200" 255 $true
$cellVal = @"

{
	"meta":
	{
		"version":"1.0.0",
		"errors":
		[
			{
				"timestamp":"2020-02-04T13:38:36.165Z
				[
					GMT
				]",
				"error":"ValidationError",
				"message":"Validation error of type FieldUndefined: Field 'geopoliticalTypeId1' in type 'GeoPoliticalType' is undefined @ 'geopoliticalTypes/geopoliticalTypeId1'",
				"path":null
			},
			{
				"timestamp":"2020-02-04T13:38:36.165Z
				[
					GMT
				]",
				"error":"ValidationError",
				"message":"Validation error of type FieldUndefined: Field 'geopoliticalTypeName2' in type 'GeoPoliticalType' is undefined @ 'geopoliticalTypes/geopoliticalTypeName2'",
				"path":null
			}
		]
	},
	"data":null
}
"@
Set-TestCell "I12" $cellVal 255 $false
Set-TestCell "J12" "Fail" 255 $false
$rng = $ws.Range("K12")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Font.Color = 255

# Reset auto-calculated row heights on the new rows back to default (no explicit height),
# matching how the source rows (2-8) have no custom height.
$ws.Range("A9:K12").EntireRow.AutoFit()

